{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove the \"Meta description: ...\" paragraph that sits right after the\n//    title (Heading1) paragraph.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Meta description\") === 0) {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 2) Insert a new bold paragraph (\"Play Banana Rock Free - A Fun and\n//    Engaging Slot Game\") right before the final (\"Prompt: ...\") paragraph.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst lastIndex = paragraphs.items.length - 1;\nconst lastPara = paragraphs.items[lastIndex];\nconst startRange = lastPara.getRange(\"Start\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +\n  '<w:t>Play Banana Rock Free - A Fun and Engaging Slot Game</w:t>' +\n  '</w:r></w:p>' +\n  '<w:p/>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nstartRange.insertOoxml(ooxml, Word.InsertLocation.before);\nawait context.sync();\n\n// The insertion above also produced a throwaway empty paragraph (needed so\n// Word treats our inserted content as a genuine, separate paragraph instead\n// of merging it into the \"Prompt: ...\" paragraph); remove that helper\n// paragraph now. It is the empty paragraph located right before the final\n// (\"Prompt: ...\") paragraph, so search backwards from the end to find it.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  if (paragraphs.items[i].text === \"\") {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 3) Replace the text of the final paragraph (previously the \"Prompt: ...\"\n//    image-generation prompt) with the meta-description copy, while keeping\n//    its existing (italic) character formatting.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst finalPara = paragraphs.items[paragraphs.items.length - 1];\nfinalPara.insertText(\n  \"Read our review of Banana Rock, the online slot game by Play'n Go that offers bonus features, free spins, and rock 'n' roll respins. Play for free now.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# 1) Remove the \"Meta description: ...\" paragraph that sits right after the\n#    title (Heading1) paragraph.\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $doc.Paragraphs.Item($i)\n    if ($p.Range.Text.IndexOf(\"Meta description\") -eq 0) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2) Insert a new bold paragraph (\"Play Banana Rock Free - A Fun and\n#    Engaging Slot Game\") right before the final (\"Prompt: ...\") paragraph.\n$count = $doc.Paragraphs.Count\n$lastPara = $doc.Paragraphs.Item($count)\n$insertionPoint = $doc.Range($lastPara.Range.Start, $lastPara.Range.Start)\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Banana Rock Free - A Fun and Engaging Slot Game</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertionPoint.InsertXML($ooxml)\n\n# The insertion above also produced a throwaway empty paragraph (needed so\n# Word treats our inserted content as a genuine, separate paragraph instead\n# of merging it into the \"Prompt: ...\" paragraph); remove that helper\n# paragraph now (its Range.Text is just the paragraph mark).\n$count = $doc.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $doc.Paragraphs.Item($i)\n    if ($p.Range.Text -eq [string][char]13) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 3) Replace the text of the final paragraph (previously the \"Prompt: ...\"\n#    image-generation prompt) with the meta-description copy, while keeping\n#    its existing (italic) character formatting.\n$count = $doc.Paragraphs.Count\n$finalPara = $doc.Paragraphs.Item($count)\n$finalRange = $finalPara.Range\n$finalRange.MoveEnd(1, -1) | Out-Null\n$finalRange.Text = \"Read our review of Banana Rock, the online slot game by Play'n Go that offers bonus features, free spins, and rock 'n' roll respins. Play for free now.\"\n"}
